$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1985294117647059
$ws.Range("C2").Value = 0.5514705882352942
$ws.Range("J2").Value = 0.01225490196078431
$ws.Range("P2").Value = 0.142156862745098
$ws.Range("S2").Value = 0.09558823529411764
$ws.Range("B3").Value = 0.008583690987124463
$ws.Range("C3").Value = 0.02145922746781116
$ws.Range("J3").Value = 0.03004291845493562
$ws.Range("P3").Value = 0.7553648068669528
$ws.Range("S3").Value = 0.184549356223176
$ws.Range("B6").Value = 0.0564516129032258
$ws.Range("D6").Value = 0.004032258064516129
$ws.Range("F6").Value = 0.04435483870967742
$ws.Range("J6").Value = 0.3024193548387097
$ws.Range("O6").Value = 0.01612903225806452
$ws.Range("Q6").Value = 0.1370967741935484
$ws.Range("R6").Value = 0.07258064516129033
$ws.Range("S6").Value = 0.3669354838709677
$ws.Range("B7").Value = 0.148
$ws.Range("D7").Value = 0.004
$ws.Range("E7").Value = 0.004
$ws.Range("F7").Value = 0.044
$ws.Range("J7").Value = 0.104
$ws.Range("O7").Value = 0.036
$ws.Range("Q7").Value = 0.124
$ws.Range("R7").Value = 0.08400000000000001
$ws.Range("S7").Value = 0.452
$ws.Range("B8").Value = 0.1244725738396625
$ws.Range("D8").Value = 0.02109704641350211
$ws.Range("F8").Value = 0.06751054852320675
$ws.Range("J8").Value = 0.1118143459915612
$ws.Range("O8").Value = 0.02742616033755274
$ws.Range("Q8").Value = 0.1624472573839662
$ws.Range("R8").Value = 0.07172995780590717
$ws.Range("S8").Value = 0.4135021097046414
$ws.Range("B9").Value = 0.1275167785234899
$ws.Range("D9").Value = 0.03355704697986577
$ws.Range("F9").Value = 0.1006711409395973
$ws.Range("J9").Value = 0.08724832214765101
$ws.Range("O9").Value = 0.02013422818791946
$ws.Range("Q9").Value = 0.1543624161073825
$ws.Range("R9").Value = 0.06711409395973154
$ws.Range("S9").Value = 0.4093959731543624
$ws.Range("B10").Value = 0.1497284716834756
$ws.Range("D10").Value = 0.01241272304111714
$ws.Range("E10").Value = 0.003103180760279286
$ws.Range("F10").Value = 0.07292474786656322
$ws.Range("J10").Value = 0.1163692785104732
$ws.Range("O10").Value = 0.01706749418153607
$ws.Range("Q10").Value = 0.1598138091543832
$ws.Range("R10").Value = 0.0791311093871218
$ws.Range("S10").Value = 0.3894491854150504
$ws.Range("G11").Value = 0.1641414141414141
$ws.Range("J11").Value = 0.08333333333333333
$ws.Range("K11").Value = 0.2272727272727273
$ws.Range("L11").Value = 0.5025252525252525
$ws.Range("S11").Value = 0.02272727272727273
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.1826923076923077
$ws.Range("L12").Value = 0.02403846153846154
$ws.Range("S12").Value = 0.04326923076923077
$ws.Range("G13").Value = 0.6739130434782609
$ws.Range("J13").Value = 0.2391304347826087
$ws.Range("S13").Value = 0.08695652173913043
$ws.Range("F15").Value = 0.02788844621513944
$ws.Range("H15").Value = 0.1354581673306773
$ws.Range("I15").Value = 0.04382470119521913
$ws.Range("J15").Value = 0.3665338645418327
$ws.Range("K15").Value = 0.08366533864541832
$ws.Range("O15").Value = 0.04780876494023904
$ws.Range("S15").Value = 0.2948207171314741
$ws.Range("F16").Value = 0.02755905511811024
$ws.Range("H16").Value = 0.1850393700787402
$ws.Range("I16").Value = 0.05511811023622047
$ws.Range("J16").Value = 0.3622047244094488
$ws.Range("K16").Value = 0.1417322834645669
$ws.Range("M16").Value = 0.02755905511811024
$ws.Range("O16").Value = 0.06299212598425197
$ws.Range("S16").Value = 0.1377952755905512
$ws.Range("F17").Value = 0.03260869565217391
$ws.Range("H17").Value = 0.1657608695652174
$ws.Range("I17").Value = 0.08423913043478261
$ws.Range("J17").Value = 0.3994565217391304
$ws.Range("K17").Value = 0.1141304347826087
$ws.Range("M17").Value = 0.0108695652173913
$ws.Range("O17").Value = 0.07065217391304347
$ws.Range("S17").Value = 0.1222826086956522
$ws.Range("F18").Value = 0.03260869565217391
$ws.Range("H18").Value = 0.1902173913043478
$ws.Range("I18").Value = 0.08152173913043478
$ws.Range("J18").Value = 0.3315217391304348
$ws.Range("K18").Value = 0.1141304347826087
$ws.Range("M18").Value = 0.03260869565217391
$ws.Range("O18").Value = 0.07608695652173914
$ws.Range("S18").Value = 0.1413043478260869
$ws.Range("F19").Value = 0.02394366197183099
$ws.Range("H19").Value = 0.2119718309859155
$ws.Range("I19").Value = 0.05704225352112676
$ws.Range("J19").Value = 0.3549295774647888
$ws.Range("K19").Value = 0.1302816901408451
$ws.Range("M19").Value = 0.02323943661971831
$ws.Range("N19").Value = 0.001408450704225352
$ws.Range("O19").Value = 0.07112676056338028
$ws.Range("S19").Value = 0.126056338028169
